$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 38464548
$ws.Range("I62").Value = 45457480
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 45457480
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -45456856
$ws.Range("N62").Value = -4648
$ws.Range("H65").Value = 38464548
$ws.Range("I65").Value = 45457480
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 227287400
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -227284280
$ws.Range("N65").Value = -23240
$ws.Range("H125").Value = 1100
$ws.Range("I125").Value = 800
$ws.Range("J125").Value = 1200
$ws.Range("K125").Value = 7200
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = -4740
$ws.Range("N125").Value = -15720
$ws.Range("H137").Value = 12502480
$ws.Range("I137").Value = 27779712
$ws.Range("K137").Value = 83339136
$ws.Range("M137").Value = -83336586
$ws.Range("H138").Value = 4373.3335
$ws.Range("I138").Value = 1449.6666
$ws.Range("J138").Value = 4999.8335
$ws.Range("K138").Value = 4348.9998
$ws.Range("L138").Value = 14999.5005
$ws.Range("M138").Value = 791.0002000000004
$ws.Range("N138").Value = -25279.5005
$ws.Range("H141").Value = 2000
$ws.Range("I141").Value = 2000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -820
$ws.Range("N141").Value = $null

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2153.2632
$ws.Range("I2").Value = 2199.9375
$ws.Range("J2").Value = 1904.3334
$ws.Range("K2").Value = 2199.9375
$ws.Range("L2").Value = 1904.3334
$ws.Range("M2").Value = -2086.9375
$ws.Range("N2").Value = -2130.3334
$ws.Range("H45").Value = 1722.6957
$ws.Range("I45").Value = 1421.8889
$ws.Range("J45").Value = 2805.6
$ws.Range("K45").Value = 1421.8889
$ws.Range("L45").Value = 2805.6
$ws.Range("M45").Value = -1044.8889
$ws.Range("N45").Value = -3559.6
$ws.Range("H61").Value = 3683.6667
$ws.Range("I61").Value = 1157.6666
$ws.Range("J61").Value = 8735.666999999999
$ws.Range("K61").Value = 1157.6666
$ws.Range("L61").Value = 8735.666999999999
$ws.Range("M61").Value = -945.6666
$ws.Range("N61").Value = -9159.666999999999
$ws.Range("H110").Value = 1008.8148
$ws.Range("I110").Value = 789.619
$ws.Range("J110").Value = 1776
$ws.Range("K110").Value = 789.619
$ws.Range("L110").Value = 1776
$ws.Range("M110").Value = 1255.381
$ws.Range("N110").Value = -5866
$ws.Range("H116").Value = 2153.2632
$ws.Range("I116").Value = 2199.9375
$ws.Range("J116").Value = 1904.3334
$ws.Range("K116").Value = 2199.9375
$ws.Range("L116").Value = 1904.3334
$ws.Range("M116").Value = 94.0625
$ws.Range("N116").Value = -6492.3334
$ws.Range("H122").Value = 1447.9131
$ws.Range("I122").Value = 1377.3636
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4132.0908
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1682.0908
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 63048.332
$ws.Range("J125").Value = 63048.332
$ws.Range("L125").Value = 63048.332
$ws.Range("N125").Value = -72888.33199999999
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060
$ws.Range("H136").Value = 3683.6667
$ws.Range("I136").Value = 1157.6666
$ws.Range("J136").Value = 8735.666999999999
$ws.Range("K136").Value = 3472.9998
$ws.Range("L136").Value = 26207.001
$ws.Range("M136").Value = -922.9998000000001
$ws.Range("N136").Value = -31307.001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2153.2632
$ws.Range("I3").Value = 2199.9375
$ws.Range("J3").Value = 1904.3334
$ws.Range("K3").Value = 2199.9375
$ws.Range("L3").Value = 1904.3334
$ws.Range("M3").Value = -2085.9375
$ws.Range("N3").Value = -2132.3334
$ws.Range("H80").Value = 124.44444
$ws.Range("I80").Value = 172.75
$ws.Range("J80").Value = 85.8
$ws.Range("K80").Value = 172.75
$ws.Range("L80").Value = 85.8
$ws.Range("M80").Value = 825.25
$ws.Range("N80").Value = -2081.8
$ws.Range("H83").Value = 124.44444
$ws.Range("I83").Value = 172.75
$ws.Range("J83").Value = 85.8
$ws.Range("K83").Value = 863.75
$ws.Range("L83").Value = 429
$ws.Range("M83").Value = 4128.25
$ws.Range("N83").Value = -10413
$ws.Range("H107").Value = 2380
$ws.Range("I107").Value = 2636.3635
$ws.Range("J107").Value = 1675
$ws.Range("K107").Value = 2636.3635
$ws.Range("L107").Value = 1675
$ws.Range("M107").Value = -716.3634999999999
$ws.Range("N107").Value = -5515
$ws.Range("H134").Value = 2706.2942
$ws.Range("I134").Value = 2006.6428
$ws.Range("J134").Value = 5971.3335
$ws.Range("K134").Value = 6019.928400000001
$ws.Range("L134").Value = 17914.0005
$ws.Range("M134").Value = -3484.928400000001
$ws.Range("N134").Value = -22984.0005

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3305.889
$ws.Range("I58").Value = 2534.1428
$ws.Range("J58").Value = 6007
$ws.Range("K58").Value = 2534.1428
$ws.Range("L58").Value = 6007
$ws.Range("M58").Value = -2331.1428
$ws.Range("N58").Value = -6413
$ws.Range("H134").Value = 2785.818
$ws.Range("I134").Value = 1650.8235
$ws.Range("K134").Value = 4952.470499999999
$ws.Range("M134").Value = -2417.470499999999
$ws.Range("H136").Value = 3305.889
$ws.Range("I136").Value = 2534.1428
$ws.Range("J136").Value = 6007
$ws.Range("K136").Value = 7602.428400000001
$ws.Range("L136").Value = 18021
$ws.Range("M136").Value = -5052.428400000001
$ws.Range("N136").Value = -23121

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1550
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 1566.6666
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 4699.9998
$ws.Range("M86").Value = -3314
$ws.Range("N86").Value = -7071.9998
$ws.Range("H89").Value = 1550
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 1566.6666
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 14099.9994
$ws.Range("M89").Value = -7572
$ws.Range("N89").Value = -25955.9994
$ws.Range("H131").Value = 2212.2603
$ws.Range("J131").Value = 2376.9849
$ws.Range("L131").Value = 7130.9547
$ws.Range("N131").Value = -17210.9547

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 16666.666
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9515
$ws.Range("H122").Value = 1666.6666
$ws.Range("I122").Value = 1666.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4999.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2549.9998
$ws.Range("N122").Value = $null
$ws.Range("H132").Value = 5255.294
$ws.Range("I132").Value = 5334.385
$ws.Range("J132").Value = 4998.25
$ws.Range("K132").Value = 16003.155
$ws.Range("L132").Value = 14994.75
$ws.Range("M132").Value = -13473.155
$ws.Range("N132").Value = -20054.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5250.3125
$ws.Range("I40").Value = 5688.8887
$ws.Range("J40").Value = 4686.4287
$ws.Range("K40").Value = 5688.8887
$ws.Range("L40").Value = 4686.4287
$ws.Range("M40").Value = -5552.8887
$ws.Range("N40").Value = -4958.4287
$ws.Range("H122").Value = 3779.5881
$ws.Range("I122").Value = 2877.5454
$ws.Range("J122").Value = 5433.3335
$ws.Range("K122").Value = 8632.636200000001
$ws.Range("L122").Value = 16300.0005
$ws.Range("M122").Value = -6182.636200000001
$ws.Range("N122").Value = -21200.0005
$ws.Range("H132").Value = 2551.6
$ws.Range("I132").Value = 2306.7302
$ws.Range("K132").Value = 6920.1906
$ws.Range("M132").Value = -4390.1906
$ws.Range("H136").Value = 2265.4285
$ws.Range("I136").Value = 1604.3529
$ws.Range("J136").Value = 5075
$ws.Range("K136").Value = 4813.0587
$ws.Range("L136").Value = 15225
$ws.Range("M136").Value = -2263.0587
$ws.Range("N136").Value = -20325

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2439.074
$ws.Range("I122").Value = 2289.6875
$ws.Range("J122").Value = 2656.3635
$ws.Range("K122").Value = 6869.0625
$ws.Range("L122").Value = 7969.0905
$ws.Range("M122").Value = -4419.0625
$ws.Range("N122").Value = -12869.0905
$ws.Range("H132").Value = 2420.5483
$ws.Range("I132").Value = 1934.5416
$ws.Range("J132").Value = 4086.8572
$ws.Range("K132").Value = 5803.6248
$ws.Range("L132").Value = 12260.5716
$ws.Range("M132").Value = -3273.6248
$ws.Range("N132").Value = -17320.5716
